$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header formatting already used by H1 (bold / bordered / centered)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data for columns I (I0) and J (IF), rows 2-38
$data = @(
    @(2, 1, 6),
    @(3, 1, 6),
    @(4, 1, 6),
    @(5, 1, 6),
    @(6, 1, 6),
    @(7, 1, 7),
    @(8, 1, 6),
    @(9, 1, 7),
    @(10, 1, 6),
    @(11, 1, 6),
    @(12, 1, 6),
    @(13, 1, 6),
    @(14, 1, 6),
    @(15, 1, 6),
    @(16, 1, 6),
    @(17, 1, 6),
    @(18, 1, 6),
    @(19, 1, 6),
    @(20, 1, 7),
    @(21, 1, 6),
    @(22, 1, 6),
    @(23, 1, 6),
    @(24, 1, 6),
    @(25, 1, 5),
    @(26, 3, 7),
    @(27, 1, 5),
    @(28, 1, 6),
    @(29, 1, 4),
    @(30, 1, 7),
    @(31, 1, 1),
    @(32, 1, 3),
    @(33, 1, 4),
    @(34, 1, 5),
    @(35, 1, 5),
    @(36, 1, 4),
    @(37, 4, 5),
    @(38, 1, 2)
)

foreach ($row in $data) {
    $r = $row[0]
    $i0 = $row[1]
    $if = $row[2]
    $ws.Cells.Item($r, 9).Value = $i0
    $ws.Cells.Item($r, 10).Value = $if
}

Write-Host "Applied I0/IF columns"
